$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values ---
$ws.Range("A1").Value = "Tom and Jerry"
$ws.Range("B1").Value = ""
$ws.Range("C1").Value = ""
$ws.Range("D1").Value = ""

$ws.Range("A2").Value = "Name"
$ws.Range("B2").Value = "Generations"
$ws.Range("C2").Value = "Generations"
$ws.Range("D2").Value = "Category"

$ws.Range("A3").Value = ""
$ws.Range("B3").Value = "Age"
$ws.Range("C3").Value = "Birth"
$ws.Range("D3").Value = ""

# --- Column widths ---
$ws.Range("A1:D3").ColumnWidth = 20.7109375

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 50
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30

# --- Merge cells ---
$ws.Range("A1:D1").Merge()
$ws.Range("A2:A3").Merge()
$ws.Range("D2:D3").Merge()

# --- Style: bold header-ish cells (A1:D1, A2:A3, C2, D2:D3, C3) ---
$headerRange = $ws.Range("A1:D3")
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.Font.Bold = $true

# --- Style: red fill + white font for B2, C2... wait only B2/B3 ---
$accent = $ws.Range("B2:B3")
$accent.Font.Bold = $false
$accent.Font.Color = 0xFFFFFF
$accent.Interior.Color = 0x0000FF

Write-Output "done"
